# Add a new "H2 balanced" scenario to the fueldata sheet.
# It is an exact copy of the "H2 heavy" scenario rows (2035 block),
# keeping every fuel/price/emission-formula value identical - only the
# Scenario label (column A) changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fueldata")

# "H2 heavy" / 2035 block currently lives in rows 34-43. Duplicate it
# into the ten rows right below (44-53), preserving styles/formulas.
$src = $ws.Range("A34:E43")
$dst = $ws.Range("A44:E53")
$src.Copy($dst)

# Re-point the VLOOKUP formulas at their own row and relabel the
# scenario as "H2 balanced" for the new block.
for ($r = 44; $r -le 53; $r++) {
    $ws.Range("A$r").Value = "H2 balanced"
    $ws.Range("E$r").Formula = '=IFERROR(VLOOKUP($C' + $r + ', data_fuelEmissions!$D$4:$F$14,3,FALSE), 0)'
}

# Leave the selection where the author last left it.
[void]$ws.Range("C50").Select()
